$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 8798.719999999999  # H43: 12725.25 -> 8798.719999999999
$ws.Cells.Item(43, 9).Value = 8810.794  # I43: 9000 -> 8810.794
$ws.Cells.Item(43, 10).Value = 8773.0625  # J43: 13967 -> 8773.0625
$ws.Cells.Item(43, 11).Value = 8810.794  # K43: 9000 -> 8810.794
$ws.Cells.Item(43, 12).Value = 8773.0625  # L43: 13967 -> 8773.0625
$ws.Cells.Item(43, 13).Value = -8741.794  # M43: -8931 -> -8741.794
$ws.Cells.Item(43, 14).Value = -8911.0625  # N43: -14105 -> -8911.0625
$ws.Cells.Item(64, 8).Value = 6173.875  # H64: 9399.6 -> 6173.875
$ws.Cells.Item(64, 9).Value = 4678.4  # I64: 0 -> 4678.4
$ws.Cells.Item(64, 10).Value = 8666.333000000001  # J64: 9399.6 -> 8666.333000000001
$ws.Cells.Item(64, 11).Value = 4678.4  # K64: 0 -> 4678.4
$ws.Cells.Item(64, 12).Value = 8666.333000000001  # L64: 9399.6 -> 8666.333000000001
$ws.Cells.Item(64, 13).Value = -4430.4  # M64: None -> -4430.4
$ws.Cells.Item(64, 14).Value = -9162.333000000001  # N64: -9895.6 -> -9162.333000000001
$ws.Cells.Item(67, 8).Value = 6173.875  # H67: 9399.6 -> 6173.875
$ws.Cells.Item(67, 9).Value = 4678.4  # I67: 0 -> 4678.4
$ws.Cells.Item(67, 10).Value = 8666.333000000001  # J67: 9399.6 -> 8666.333000000001
$ws.Cells.Item(67, 11).Value = 4678.4  # K67: 0 -> 4678.4
$ws.Cells.Item(67, 12).Value = 8666.333000000001  # L67: 9399.6 -> 8666.333000000001
$ws.Cells.Item(67, 13).Value = -3820.4  # M67: None -> -3820.4
$ws.Cells.Item(67, 14).Value = -10382.333  # N67: -11115.6 -> -10382.333
$ws.Cells.Item(106, 8).Value = 6311.1113  # H106: 6311.778 -> 6311.1113
$ws.Cells.Item(106, 9).Value = 7167.1665  # I106: 7168.1665 -> 7167.1665
$ws.Cells.Item(106, 11).Value = 7167.1665  # K106: 7168.1665 -> 7167.1665
$ws.Cells.Item(106, 13).Value = -6536.1665  # M106: -6537.1665 -> -6536.1665
$ws.Cells.Item(116, 8).Value = 3159.9  # H116: 3118.875 -> 3159.9
$ws.Cells.Item(116, 9).Value = 2915.2632  # I116: 2915.4092 -> 2915.2632
$ws.Cells.Item(116, 10).Value = 7808  # J116: 5357 -> 7808
$ws.Cells.Item(116, 11).Value = 2915.2632  # K116: 2915.4092 -> 2915.2632
$ws.Cells.Item(116, 12).Value = 7808  # L116: 5357 -> 7808
$ws.Cells.Item(116, 13).Value = 526.7368000000001  # M116: 526.5907999999999 -> 526.7368000000001
$ws.Cells.Item(116, 14).Value = -14692  # N116: -12241 -> -14692
$ws.Cells.Item(132, 8).Value = 13392.958  # H132: 14680.5 -> 13392.958
$ws.Cells.Item(132, 9).Value = 10726.85  # I132: 11948.056 -> 10726.85
$ws.Cells.Item(132, 10).Value = 26723.5  # J132: 26976.5 -> 26723.5
$ws.Cells.Item(132, 11).Value = 32180.55  # K132: 35844.16800000001 -> 32180.55
$ws.Cells.Item(132, 12).Value = 80170.5  # L132: 80929.5 -> 80170.5
$ws.Cells.Item(132, 13).Value = -29650.55  # M132: -33314.16800000001 -> -29650.55
$ws.Cells.Item(132, 14).Value = -85230.5  # N132: -85989.5 -> -85230.5
$ws.Cells.Item(135, 8).Value = 1907.575  # H135: 1908.925 -> 1907.575
$ws.Cells.Item(135, 9).Value = 617.63635  # I135: 619.2727 -> 617.63635
$ws.Cells.Item(135, 11).Value = 5558.72715  # K135: 5573.454299999999 -> 5558.72715
$ws.Cells.Item(135, 13).Value = -3023.72715  # M135: -3038.454299999999 -> -3023.72715
$ws.Cells.Item(137, 8).Value = 6498.7754  # H137: 6922.391 -> 6498.7754
$ws.Cells.Item(137, 9).Value = 2364.5186  # I137: 2574.2 -> 2364.5186
$ws.Cells.Item(137, 10).Value = 11572.637  # J137: 12098.81 -> 11572.637
$ws.Cells.Item(137, 11).Value = 7093.5558  # K137: 7722.599999999999 -> 7093.5558
$ws.Cells.Item(137, 12).Value = 34717.911  # L137: 36296.43 -> 34717.911
$ws.Cells.Item(137, 13).Value = -4543.5558  # M137: -5172.599999999999 -> -4543.5558
$ws.Cells.Item(137, 14).Value = -39817.911  # N137: -41396.43 -> -39817.911
$ws.Cells.Item(138, 8).Value = 2415.21  # H138: 2397.21 -> 2415.21
$ws.Cells.Item(138, 10).Value = 2670.3613  # J138: 2648.6748 -> 2670.3613
$ws.Cells.Item(138, 12).Value = 8011.0839  # L138: 7946.024399999999 -> 8011.0839
$ws.Cells.Item(138, 14).Value = -18291.0839  # N138: -18226.0244 -> -18291.0839

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 880.5  # H12: 687.375 -> 880.5
$ws.Cells.Item(12, 9).Value = 779  # I12: 470 -> 779
$ws.Cells.Item(12, 11).Value = 779  # K12: 470 -> 779
$ws.Cells.Item(12, 13).Value = -606  # M12: -297 -> -606
$ws.Cells.Item(32, 8).Value = 6757.17  # H32: 6608.19 -> 6757.17
$ws.Cells.Item(32, 9).Value = 1243.5256  # I32: 1244.8334 -> 1243.5256
$ws.Cells.Item(32, 10).Value = 26305.545  # J32: 25623.727 -> 26305.545
$ws.Cells.Item(32, 11).Value = 1243.5256  # K32: 1244.8334 -> 1243.5256
$ws.Cells.Item(32, 12).Value = 26305.545  # L32: 25623.727 -> 26305.545
$ws.Cells.Item(32, 13).Value = -956.5255999999999  # M32: -957.8334 -> -956.5255999999999
$ws.Cells.Item(32, 14).Value = -26879.545  # N32: -26197.727 -> -26879.545
$ws.Cells.Item(61, 8).Value = 11832.108  # H61: 13940.719 -> 11832.108
$ws.Cells.Item(61, 9).Value = 1584.9546  # I61: 2029 -> 1584.9546
$ws.Cells.Item(61, 10).Value = 26861.268  # J61: 23205.389 -> 26861.268
$ws.Cells.Item(61, 11).Value = 1584.9546  # K61: 2029 -> 1584.9546
$ws.Cells.Item(61, 12).Value = 26861.268  # L61: 23205.389 -> 26861.268
$ws.Cells.Item(61, 13).Value = -1372.9546  # M61: -1817 -> -1372.9546
$ws.Cells.Item(61, 14).Value = -27285.268  # N61: -23629.389 -> -27285.268
$ws.Cells.Item(136, 8).Value = 11832.108  # H136: 13940.719 -> 11832.108
$ws.Cells.Item(136, 9).Value = 1584.9546  # I136: 2029 -> 1584.9546
$ws.Cells.Item(136, 10).Value = 26861.268  # J136: 23205.389 -> 26861.268
$ws.Cells.Item(136, 11).Value = 4754.8638  # K136: 6087 -> 4754.8638
$ws.Cells.Item(136, 12).Value = 80583.804  # L136: 69616.167 -> 80583.804
$ws.Cells.Item(136, 13).Value = -2204.8638  # M136: -3537 -> -2204.8638
$ws.Cells.Item(136, 14).Value = -85683.804  # N136: -74716.167 -> -85683.804

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 5886.65  # H86: 76932960 -> 5886.65
$ws.Cells.Item(86, 9).Value = 6959.3125  # I86: 12448.7 -> 6959.3125
$ws.Cells.Item(86, 10).Value = 1596  # J86: 333334660 -> 1596
$ws.Cells.Item(86, 11).Value = 6959.3125  # K86: 12448.7 -> 6959.3125
$ws.Cells.Item(86, 12).Value = 1596  # L86: 333334660 -> 1596
$ws.Cells.Item(86, 13).Value = -5836.3125  # M86: -11325.7 -> -5836.3125
$ws.Cells.Item(86, 14).Value = -3842  # N86: -333336906 -> -3842
$ws.Cells.Item(89, 8).Value = 5886.65  # H89: 76932960 -> 5886.65
$ws.Cells.Item(89, 9).Value = 6959.3125  # I89: 12448.7 -> 6959.3125
$ws.Cells.Item(89, 10).Value = 1596  # J89: 333334660 -> 1596
$ws.Cells.Item(89, 11).Value = 34796.5625  # K89: 62243.5 -> 34796.5625
$ws.Cells.Item(89, 12).Value = 7980  # L89: 1666673300 -> 7980
$ws.Cells.Item(89, 13).Value = -29180.5625  # M89: -56627.5 -> -29180.5625
$ws.Cells.Item(89, 14).Value = -19212  # N89: -1666684532 -> -19212
$ws.Cells.Item(94, 8).Value = 3295.8262  # H94: 2753.4092 -> 3295.8262
$ws.Cells.Item(94, 9).Value = 2827.8572  # I94: 1655.5294 -> 2827.8572
$ws.Cells.Item(94, 10).Value = 4023.7778  # J94: 6486.2 -> 4023.7778
$ws.Cells.Item(94, 11).Value = 2827.8572  # K94: 1655.5294 -> 2827.8572
$ws.Cells.Item(94, 12).Value = 4023.7778  # L94: 6486.2 -> 4023.7778
$ws.Cells.Item(94, 13).Value = -2376.8572  # M94: -1204.5294 -> -2376.8572
$ws.Cells.Item(94, 14).Value = -4925.7778  # N94: -7388.2 -> -4925.7778
$ws.Cells.Item(133, 8).Value = 95777.5  # H133: 122890 -> 95777.5
$ws.Cells.Item(133, 10).Value = 95777.5  # J133: 122890 -> 95777.5
$ws.Cells.Item(133, 12).Value = 95777.5  # L133: 122890 -> 95777.5
$ws.Cells.Item(133, 14).Value = -105897.5  # N133: -133010 -> -105897.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(11, 8).Value = 0  # H11: 307.7 -> 0
$ws.Cells.Item(11, 9).Value = 0  # I11: 304.80768 -> 0
$ws.Cells.Item(11, 10).Value = 0  # J11: 326.5 -> 0
$ws.Cells.Item(11, 11).Value = 0  # K11: 304.80768 -> 0
$ws.Cells.Item(11, 12).Value = 0  # L11: 326.5 -> 0
$ws.Cells.Item(11, 13).ClearContents()  # M11: was -164.80768
$ws.Cells.Item(11, 14).ClearContents()  # N11: was -606.5
$ws.Cells.Item(12, 8).Value = 346.6875  # H12: 1077.6 -> 346.6875
$ws.Cells.Item(12, 9).Value = 52.076923  # I12: 378.8 -> 52.076923
$ws.Cells.Item(12, 10).Value = 1623.3334  # J12: 3174 -> 1623.3334
$ws.Cells.Item(12, 11).Value = 52.076923  # K12: 378.8 -> 52.076923
$ws.Cells.Item(12, 12).Value = 1623.3334  # L12: 3174 -> 1623.3334
$ws.Cells.Item(12, 13).Value = 117.923077  # M12: -208.8 -> 117.923077
$ws.Cells.Item(12, 14).Value = -1963.3334  # N12: -3514 -> -1963.3334
$ws.Cells.Item(14, 8).Value = 0  # H14: 400 -> 0
$ws.Cells.Item(14, 9).Value = 0  # I14: 400 -> 0
$ws.Cells.Item(14, 10).Value = 0  # J14: 400 -> 0
$ws.Cells.Item(14, 11).Value = 0  # K14: 400 -> 0
$ws.Cells.Item(14, 12).Value = 0  # L14: 400 -> 0
$ws.Cells.Item(14, 13).ClearContents()  # M14: was -230
$ws.Cells.Item(14, 14).ClearContents()  # N14: was -740
$ws.Cells.Item(22, 8).Value = 734.875  # H22: 953.36365 -> 734.875
$ws.Cells.Item(22, 9).Value = 415  # I22: 436.25 -> 415
$ws.Cells.Item(22, 10).Value = 1694.5  # J22: 2332.3333 -> 1694.5
$ws.Cells.Item(22, 11).Value = 415  # K22: 436.25 -> 415
$ws.Cells.Item(22, 12).Value = 1694.5  # L22: 2332.3333 -> 1694.5
$ws.Cells.Item(22, 13).Value = -65  # M22: -86.25 -> -65
$ws.Cells.Item(22, 14).Value = -2394.5  # N22: -3032.3333 -> -2394.5
$ws.Cells.Item(23, 8).Value = 10250  # H23: 338383.34 -> 10250
$ws.Cells.Item(23, 9).Value = 5000  # I23: 1000000 -> 5000
$ws.Cells.Item(23, 10).Value = 15500  # J23: 7575 -> 15500
$ws.Cells.Item(23, 11).Value = 5000  # K23: 1000000 -> 5000
$ws.Cells.Item(23, 12).Value = 15500  # L23: 7575 -> 15500
$ws.Cells.Item(23, 13).Value = -4760  # M23: -999760 -> -4760
$ws.Cells.Item(23, 14).Value = -15980  # N23: -8055 -> -15980
$ws.Cells.Item(27, 8).Value = 10250  # H27: 338383.34 -> 10250
$ws.Cells.Item(27, 9).Value = 5000  # I27: 1000000 -> 5000
$ws.Cells.Item(27, 10).Value = 15500  # J27: 7575 -> 15500
$ws.Cells.Item(27, 11).Value = 5000  # K27: 1000000 -> 5000
$ws.Cells.Item(27, 12).Value = 15500  # L27: 7575 -> 15500
$ws.Cells.Item(27, 13).Value = -4808  # M27: -999808 -> -4808
$ws.Cells.Item(27, 14).Value = -15884  # N27: -7959 -> -15884
$ws.Cells.Item(31, 8).Value = 10088.633  # H31: 9479.537 -> 10088.633
$ws.Cells.Item(31, 9).Value = 4085.1282  # I31: 3946.9285 -> 4085.1282
$ws.Cells.Item(31, 10).Value = 33502.3  # J31: 28843.666 -> 33502.3
$ws.Cells.Item(31, 11).Value = 4085.1282  # K31: 3946.9285 -> 4085.1282
$ws.Cells.Item(31, 12).Value = 33502.3  # L31: 28843.666 -> 33502.3
$ws.Cells.Item(31, 13).Value = -3790.1282  # M31: -3651.9285 -> -3790.1282
$ws.Cells.Item(31, 14).Value = -34092.3  # N31: -29433.666 -> -34092.3
$ws.Cells.Item(34, 8).Value = 10088.633  # H34: 9479.537 -> 10088.633
$ws.Cells.Item(34, 9).Value = 4085.1282  # I34: 3946.9285 -> 4085.1282
$ws.Cells.Item(34, 10).Value = 33502.3  # J34: 28843.666 -> 33502.3
$ws.Cells.Item(34, 11).Value = 4085.1282  # K34: 3946.9285 -> 4085.1282
$ws.Cells.Item(34, 12).Value = 33502.3  # L34: 28843.666 -> 33502.3
$ws.Cells.Item(34, 13).Value = -3883.1282  # M34: -3744.9285 -> -3883.1282
$ws.Cells.Item(34, 14).Value = -33906.3  # N34: -29247.666 -> -33906.3
$ws.Cells.Item(38, 8).Value = 3622  # H38: 100000 -> 3622
$ws.Cells.Item(38, 9).Value = 3622  # I38: 100000 -> 3622
$ws.Cells.Item(38, 11).Value = 3622  # K38: 100000 -> 3622
$ws.Cells.Item(38, 13).Value = -3245  # M38: -99623 -> -3245
$ws.Cells.Item(41, 8).Value = 34999  # H41: 23933.334 -> 34999
$ws.Cells.Item(41, 9).Value = 34999  # I41: 23933.334 -> 34999
$ws.Cells.Item(41, 11).Value = 34999  # K41: 23933.334 -> 34999
$ws.Cells.Item(41, 13).Value = -34571  # M41: -23505.334 -> -34571
$ws.Cells.Item(46, 8).Value = 3622  # H46: 100000 -> 3622
$ws.Cells.Item(46, 9).Value = 3622  # I46: 100000 -> 3622
$ws.Cells.Item(46, 11).Value = 3622  # K46: 100000 -> 3622
$ws.Cells.Item(46, 13).Value = -3411  # M46: -99789 -> -3411
$ws.Cells.Item(58, 8).Value = 10387.234  # H58: 10321.83 -> 10387.234
$ws.Cells.Item(58, 9).Value = 3935.3928  # I58: 3929.0715 -> 3935.3928
$ws.Cells.Item(58, 10).Value = 19895.21  # J58: 19742.736 -> 19895.21
$ws.Cells.Item(58, 11).Value = 3935.3928  # K58: 3929.0715 -> 3935.3928
$ws.Cells.Item(58, 12).Value = 19895.21  # L58: 19742.736 -> 19895.21
$ws.Cells.Item(58, 13).Value = -3732.3928  # M58: -3726.0715 -> -3732.3928
$ws.Cells.Item(58, 14).Value = -20301.21  # N58: -20148.736 -> -20301.21
$ws.Cells.Item(99, 8).Value = 7473.6313  # H99: 5595.9644 -> 7473.6313
$ws.Cells.Item(99, 9).Value = 1999  # I99: 1576.5555 -> 1999
$ws.Cells.Item(99, 10).Value = 7777.778  # J99: 7499.8945 -> 7777.778
$ws.Cells.Item(99, 11).Value = 1999  # K99: 1576.5555 -> 1999
$ws.Cells.Item(99, 12).Value = 7777.778  # L99: 7499.8945 -> 7777.778
$ws.Cells.Item(99, 13).Value = -501  # M99: -78.55549999999994 -> -501
$ws.Cells.Item(99, 14).Value = -10773.778  # N99: -10495.8945 -> -10773.778
$ws.Cells.Item(126, 8).Value = 7473.6313  # H126: 5595.9644 -> 7473.6313
$ws.Cells.Item(126, 9).Value = 1999  # I126: 1576.5555 -> 1999
$ws.Cells.Item(126, 10).Value = 7777.778  # J126: 7499.8945 -> 7777.778
$ws.Cells.Item(126, 11).Value = 5997  # K126: 4729.666499999999 -> 5997
$ws.Cells.Item(126, 12).Value = 23333.334  # L126: 22499.6835 -> 23333.334
$ws.Cells.Item(126, 13).Value = -3527  # M126: -2259.666499999999 -> -3527
$ws.Cells.Item(126, 14).Value = -28273.334  # N126: -27439.6835 -> -28273.334
$ws.Cells.Item(132, 8).Value = 7542.125  # H132: 8861.4 -> 7542.125
$ws.Cells.Item(132, 9).Value = 2681.6365  # I132: 3600.1428 -> 2681.6365
$ws.Cells.Item(132, 10).Value = 11654.846  # J132: 11694.385 -> 11654.846
$ws.Cells.Item(132, 11).Value = 8044.9095  # K132: 10800.4284 -> 8044.9095
$ws.Cells.Item(132, 12).Value = 34964.538  # L132: 35083.155 -> 34964.538
$ws.Cells.Item(132, 13).Value = -5514.9095  # M132: -8270.428400000001 -> -5514.9095
$ws.Cells.Item(132, 14).Value = -40024.538  # N132: -40143.155 -> -40024.538
$ws.Cells.Item(134, 8).Value = 22732524  # H134: 18523114 -> 22732524
$ws.Cells.Item(134, 9).Value = 1512.9231  # I134: 1313.3529 -> 1512.9231
$ws.Cells.Item(134, 10).Value = 55566204  # J134: 50010172 -> 55566204
$ws.Cells.Item(134, 11).Value = 4538.7693  # K134: 3940.0587 -> 4538.7693
$ws.Cells.Item(134, 12).Value = 166698612  # L134: 150030516 -> 166698612
$ws.Cells.Item(134, 13).Value = -2003.7693  # M134: -1405.0587 -> -2003.7693
$ws.Cells.Item(134, 14).Value = -166703682  # N134: -150035586 -> -166703682
$ws.Cells.Item(136, 8).Value = 10387.234  # H136: 10321.83 -> 10387.234
$ws.Cells.Item(136, 9).Value = 3935.3928  # I136: 3929.0715 -> 3935.3928
$ws.Cells.Item(136, 10).Value = 19895.21  # J136: 19742.736 -> 19895.21
$ws.Cells.Item(136, 11).Value = 11806.1784  # K136: 11787.2145 -> 11806.1784
$ws.Cells.Item(136, 12).Value = 59685.63  # L136: 59228.208 -> 59685.63
$ws.Cells.Item(136, 13).Value = -9256.178400000001  # M136: -9237.2145 -> -9256.178400000001
$ws.Cells.Item(136, 14).Value = -64785.63  # N136: -64328.208 -> -64785.63

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 560.7692  # H23: 543.55554 -> 560.7692
$ws.Cells.Item(23, 9).Value = 615  # I23: 598.2727 -> 615
$ws.Cells.Item(23, 10).Value = 526.875  # J23: 505.9375 -> 526.875
$ws.Cells.Item(23, 11).Value = 1845  # K23: 1794.8181 -> 1845
$ws.Cells.Item(23, 12).Value = 1580.625  # L23: 1517.8125 -> 1580.625
$ws.Cells.Item(23, 13).Value = -1610  # M23: -1559.8181 -> -1610
$ws.Cells.Item(23, 14).Value = -2050.625  # N23: -1987.8125 -> -2050.625
$ws.Cells.Item(113, 8).Value = 1306.1765  # H113: 1376.5 -> 1306.1765
$ws.Cells.Item(113, 9).Value = 620.6  # I113: 620.4 -> 620.6
$ws.Cells.Item(113, 10).Value = 1591.8334  # J113: 1796.5555 -> 1591.8334
$ws.Cells.Item(113, 11).Value = 1861.8  # K113: 1861.2 -> 1861.8
$ws.Cells.Item(113, 12).Value = 4775.5002  # L113: 5389.666499999999 -> 4775.5002
$ws.Cells.Item(113, 13).Value = 308.1999999999998  # M113: 308.8000000000002 -> 308.1999999999998
$ws.Cells.Item(113, 14).Value = -9115.5002  # N113: -9729.666499999999 -> -9115.5002
$ws.Cells.Item(136, 8).Value = 2000  # H136: 2197.6 -> 2000
$ws.Cells.Item(136, 9).Value = 2000  # I136: 1772 -> 2000
$ws.Cells.Item(136, 10).Value = 0  # J136: 3900 -> 0
$ws.Cells.Item(136, 11).Value = 6000  # K136: 5316 -> 6000
$ws.Cells.Item(136, 12).Value = 0  # L136: 11700 -> 0
$ws.Cells.Item(136, 13).Value = -900  # M136: -216 -> -900
$ws.Cells.Item(136, 14).ClearContents()  # N136: was -21900

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 100  # H12: 855.38464 -> 100
$ws.Cells.Item(12, 9).Value = 100  # I12: 93.333336 -> 100
$ws.Cells.Item(12, 10).Value = 0  # J12: 10000 -> 0
$ws.Cells.Item(12, 11).Value = 100  # K12: 93.333336 -> 100
$ws.Cells.Item(12, 12).Value = 0  # L12: 10000 -> 0
$ws.Cells.Item(12, 13).Value = 40  # M12: 46.666664 -> 40
$ws.Cells.Item(12, 14).ClearContents()  # N12: was -10280

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 17115.69  # H136: 20417.084 -> 17115.69
$ws.Cells.Item(136, 9).Value = 20653.455  # I136: 43078.6 -> 20653.455
$ws.Cells.Item(136, 10).Value = 14953.723  # J136: 14453.526 -> 14953.723
$ws.Cells.Item(136, 11).Value = 61960.36500000001  # K136: 129235.8 -> 61960.36500000001
$ws.Cells.Item(136, 12).Value = 44861.169  # L136: 43360.578 -> 44861.169
$ws.Cells.Item(136, 13).Value = -59410.36500000001  # M136: -126685.8 -> -59410.36500000001
$ws.Cells.Item(136, 14).Value = -49961.169  # N136: -48460.578 -> -49961.169

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2690.8  # H96: 2487.5 -> 2690.8
$ws.Cells.Item(96, 9).Value = 1500  # I96: 0 -> 1500
$ws.Cells.Item(96, 10).Value = 2988.5  # J96: 2487.5 -> 2988.5
$ws.Cells.Item(96, 11).Value = 1500  # K96: 0 -> 1500
$ws.Cells.Item(96, 12).Value = 2988.5  # L96: 2487.5 -> 2988.5
$ws.Cells.Item(96, 13).Value = -127  # M96: None -> -127
$ws.Cells.Item(96, 14).Value = -5734.5  # N96: -5233.5 -> -5734.5
$ws.Cells.Item(113, 8).Value = 790.9655  # H113: 2819.8845 -> 790.9655
$ws.Cells.Item(113, 9).Value = 489.47058  # I113: 3136.0527 -> 489.47058
$ws.Cells.Item(113, 10).Value = 1218.0834  # J113: 1961.7142 -> 1218.0834
$ws.Cells.Item(113, 11).Value = 1468.41174  # K113: 9408.158100000001 -> 1468.41174
$ws.Cells.Item(113, 12).Value = 3654.2502  # L113: 5885.142599999999 -> 3654.2502
$ws.Cells.Item(113, 13).Value = 701.58826  # M113: -7238.158100000001 -> 701.58826
$ws.Cells.Item(113, 14).Value = -7994.2502  # N113: -10225.1426 -> -7994.2502
$ws.Cells.Item(122, 8).Value = 2964.125  # H122: 2258.8823 -> 2964.125
$ws.Cells.Item(122, 9).Value = 1129.6  # I122: 960.28 -> 1129.6
$ws.Cells.Item(122, 10).Value = 6021.6665  # J122: 5866.1113 -> 6021.6665
$ws.Cells.Item(122, 11).Value = 3388.8  # K122: 2880.84 -> 3388.8
$ws.Cells.Item(122, 12).Value = 18064.9995  # L122: 17598.3339 -> 18064.9995
$ws.Cells.Item(122, 13).Value = -938.7999999999997  # M122: -430.8400000000001 -> -938.7999999999997
$ws.Cells.Item(122, 14).Value = -22964.9995  # N122: -22498.3339 -> -22964.9995
$ws.Cells.Item(132, 8).Value = 6598.9375  # H132: 5739.3115 -> 6598.9375
$ws.Cells.Item(132, 9).Value = 3537.6924  # I132: 2924.2974 -> 3537.6924
$ws.Cells.Item(132, 10).Value = 10216.772  # J132: 10079.125 -> 10216.772
$ws.Cells.Item(132, 11).Value = 10613.0772  # K132: 8772.8922 -> 10613.0772
$ws.Cells.Item(132, 12).Value = 30650.316  # L132: 30237.375 -> 30650.316
$ws.Cells.Item(132, 13).Value = -8083.0772  # M132: -6242.8922 -> -8083.0772
$ws.Cells.Item(132, 14).Value = -35710.31600000001  # N132: -35297.375 -> -35710.31600000001
